# Update employee absence data rows 2-11 with new values as per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 24731
$ws.Cells.Item(2, 2).Value = "Beatriz Souza"
$ws.Cells.Item(2, 3).Value = "TI"
$ws.Cells.Item(2, 4).Value = "Outros"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 45089
$ws.Cells.Item(2, 7).Value = 5225.59

# Row 3
$ws.Cells.Item(3, 1).Value = 14764
$ws.Cells.Item(3, 2).Value = "Zoe Leão"
$ws.Cells.Item(3, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(3, 4).Value = "Doenca"
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 45091
$ws.Cells.Item(3, 7).Value = 7388.23

# Row 4
$ws.Cells.Item(4, 1).Value = 6605
$ws.Cells.Item(4, 2).Value = "Letícia Borges"
$ws.Cells.Item(4, 3).Value = "P&D"
$ws.Cells.Item(4, 4).Value = "Outros"
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 45092
$ws.Cells.Item(4, 7).Value = 6471.21

# Row 5
$ws.Cells.Item(5, 1).Value = 85721
$ws.Cells.Item(5, 2).Value = "Anna Liz Siqueira"
$ws.Cells.Item(5, 3).Value = "Engenharia"
$ws.Cells.Item(5, 4).Value = "Outros"
$ws.Cells.Item(5, 5).Value = 6
$ws.Cells.Item(5, 6).Value = 45103
$ws.Cells.Item(5, 7).Value = 2645.2

# Row 6
$ws.Cells.Item(6, 1).Value = 70216
$ws.Cells.Item(6, 2).Value = "Maria Júlia da Mota"
$ws.Cells.Item(6, 3).Value = "P&D"
$ws.Cells.Item(6, 4).Value = "Problemas pessoais"
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 45084
$ws.Cells.Item(6, 7).Value = 5348.62

# Row 7
$ws.Cells.Item(7, 1).Value = 80740
$ws.Cells.Item(7, 2).Value = "Melina Araújo"
$ws.Cells.Item(7, 3).Value = "Financeiro"
$ws.Cells.Item(7, 4).Value = "Doenca"
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 45083
$ws.Cells.Item(7, 7).Value = 6682.61

# Row 8
$ws.Cells.Item(8, 1).Value = 94788
$ws.Cells.Item(8, 2).Value = "Dante da Paz"
$ws.Cells.Item(8, 3).Value = "Juridico"
$ws.Cells.Item(8, 4).Value = "Problemas pessoais"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 45083
$ws.Cells.Item(8, 7).Value = 4598.1

# Row 9
$ws.Cells.Item(9, 1).Value = 53887
$ws.Cells.Item(9, 2).Value = "Thiago Pastor"
$ws.Cells.Item(9, 3).Value = "Recursos Humanos"
$ws.Cells.Item(9, 4).Value = "Problemas pessoais"
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 45078
$ws.Cells.Item(9, 7).Value = 5080.86

# Row 10
$ws.Cells.Item(10, 1).Value = 48074
$ws.Cells.Item(10, 2).Value = "Marcos Vinicius da Mata"
$ws.Cells.Item(10, 3).Value = "P&D"
$ws.Cells.Item(10, 4).Value = "Problemas pessoais"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 45091
$ws.Cells.Item(10, 7).Value = 7853.4

# Row 11
$ws.Cells.Item(11, 1).Value = 28371
$ws.Cells.Item(11, 2).Value = "Milena Camargo"
$ws.Cells.Item(11, 3).Value = "Financeiro"
$ws.Cells.Item(11, 4).Value = "Consulta medica"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 45089
$ws.Cells.Item(11, 7).Value = 2862.47
